$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new consecutive invoice number (zero-padded) and updated entry/exit times
$ws.Range("A2").Value = "FE-0317055"
$ws.Range("C2").Value = "29/07/2024 15:57:00"
$ws.Range("D2").Value = "29/07/2024 15:58:00"
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 1500

# Row 3: new invoice number, plate, updated entry/exit times and elapsed time
$ws.Range("A3").Value = "FE-0317056"
$ws.Range("B3").Value = "BBB01"
$ws.Range("C3").Value = "29/07/2024 16:38:00"
$ws.Range("D3").Value = "29/07/2024 16:55:00"
$ws.Range("G3").Value = 1020

# Row 4: new invoice number, plate, updated entry/exit times, value, elapsed time and total
$ws.Range("A4").Value = "FE-0317057"
$ws.Range("B4").Value = "CCC01"
$ws.Range("C4").Value = "29/07/2024 17:10:00"
$ws.Range("D4").Value = "29/07/2024 18:27:00"
$ws.Range("F4").Value = 1500
$ws.Range("G4").Value = 4620
$ws.Range("H4").Value = 3000

# Rows 5 and 6 no longer exist in the updated register — remove them entirely
$ws.Range("A5:H6").EntireRow.Delete()
